# excel update 2020-04-03 10:30 + added demo & hospitalization
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Status" sheet: append the 2020-04-03 row, extend the trailing blank
#    formatted rows, and update the active selection.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Status")

$status.Range("A59").Value = 43924.4375
$status.Range("B59").Value = 62408
$status.Range("C59").Value = 1245
$status.Range("D59").Value = 3255
$status.Range("D59").Style = "Normal"
$status.Range("E59").Value = 1023
$status.Range("F59").Value = 67
$status.Range("G59").Value = 67998
$status.Range("G59").Style = "Normal"

# Row 64 gains a number-formatted (but still empty) D cell, matching the
# existing F64/G64 "#,##0" formatting.
$status.Range("D64").NumberFormat = "#,##0"

# A brand-new blank row 66 with the same single formatted D cell as its
# neighbours (rows 65/67).
$status.Range("D66").NumberFormat = "#,##0"

$status.Activate()
$status.Range("A59").Select()

# ---------------------------------------------------------------------------
# 2. Insert the new "Demographics" sheet right after "Status" (and before
#    "Daily Summary"), then hide "Daily Summary".
# ---------------------------------------------------------------------------
$demo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $status)
$demo.Name = "Demographics"

$dailySummary = $wb.Worksheets.Item("Daily Summary")
$dailySummary.Visible = $false

# ---------------------------------------------------------------------------
# 3. Populate "Demographics".
# ---------------------------------------------------------------------------

# -- Header row -------------------------------------------------------------
$demo.Range("A1").Value = "Date"
$demo.Range("A1").Font.Name = "Calibri"
$demo.Range("A1").Font.Size = 12
$demo.Range("A1").Font.Color = 0

$demo.Range("B1").Value = "Male"
$demo.Range("C1").Value = "Female"
$demo.Range("D1").Value = "19 and under"
$demo.Range("E1").Value = "20-39"
$demo.Range("F1").Value = "40-59"
$demo.Range("G1").Value = "60-79"
$demo.Range("H1").Value = "80 and over"
$demo.Range("I1").Value = "Number of patients hospitalized with COVID-196"
$demo.Range("J1").Value = "Number of patients in ICU7 with COVID-19"
$demo.Range("K1").Value = "Number of patients in ICU7 on a ventilator with COVID-19"

# -- Data rows (2 and 3) -----------------------------------------------------
# Seed the values first, then copy the existing date / thousands-style
# formats from the "Status" sheet so the new cells re-use the workbook's
# existing style entries instead of minting new ones.
$demo.Range("A2").Value = 43923.4375
$demo.Range("A3").Value = 43924.4375
$status.Range("A2").Copy()
$demo.Range("A2:A3").PasteSpecial(-4122)

$demo.Range("B2").Value = 1355
$demo.Range("C2").Value = 1421
$demo.Range("D2").Value = 68
$demo.Range("E2").Value = 826
$demo.Range("F2").Value = 995
$demo.Range("G2").Value = 716
$demo.Range("H2").Value = 186

$demo.Range("B3").Value = 1579
$demo.Range("C3").Value = 1657
$demo.Range("D3").Value = 82
$demo.Range("E3").Value = 945
$demo.Range("F3").Value = 1178
$demo.Range("G3").Value = 821
$demo.Range("H3").Value = 226
$demo.Range("I3").Value = 462
$demo.Range("J3").Value = 194
$demo.Range("K3").Value = 140

$status.Range("D60").Copy()
$demo.Range("B3").PasteSpecial(-4122)
$demo.Range("C3").PasteSpecial(-4122)
$demo.Range("F3").PasteSpecial(-4122)

# -- Five "card" rows (4-8) ---------------------------------------------------
$demo.Range("A4").Font.Name = "Helvetica Neue"
$demo.Range("A4").Font.Size = 16
$demo.Range("A4").Font.Color = 2236962
$demo.Range("A4").Copy()
$demo.Range("A4:C8").PasteSpecial(-4122)
for ($r = 4; $r -le 8; $r++) {
    $demo.Rows.Item($r).RowHeight = 20
}

# -- Scattered "#,##0"-formatted placeholder cells ---------------------------
$status.Range("D60").Copy()
$demo.Range("G9").PasteSpecial(-4122)
$demo.Range("G10").PasteSpecial(-4122)
$demo.Range("G13").PasteSpecial(-4122)
$demo.Range("G17").PasteSpecial(-4122)
$demo.Range("G18").PasteSpecial(-4122)

# -- Column width / view settings -------------------------------------------
$demo.Columns.Item(1).ColumnWidth = 14.6666666666667

$demo.Activate()
$demo.Range("F8:H22").Select()
$excel.ActiveWindow.Zoom = 150
